# FA170_TestData_CreateAccounting_21C.xlsx - re-upload edit
# Clears the stored Oracle Cloud login/URL values on the Input_Value sheet
# (URL / UserName / Password) before re-sharing the workbook, and leaves the
# sheet scrolled/selected near those cells, matching how the sheet was left
# when it was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

$ws.Range("S2").Value = ""
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = ""

$ws.Range("S2:U2").Select()
